# Update JSON range function and host files
#
# Re-creates the commit that:
#  - moves the cursor selection,
#  - collapses the A5:L6 "_xldudf_STREAM_PARSEJSONRANGE(A3)" spill down to a
#    single row (A5:L5) now that the upstream JSON only has one record, and
#    removes the now-unused row 6,
#  - adds a second TRUE argument to the A10 array formula
#    "_xldudf_STREAM_PARSEJSONRANGE(_xldudf_STREAM_EXTRACT(A1,"data"),TRUE)".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 5 / Row 6 -------------------------------------------------------
# Target layout: the array formula in A5 now only spills across A5:L5 (the
# single ACE.AX record) instead of A5:L6, and row 6 (the old second record)
# disappears.

# Write the literal "spilled" values for B5:L5 first (these are safe plain
# values, not covered by any array yet) so they survive the re-entry of the
# array formula below.
$ws.Range("B5").Value2 = "Acusensus Ltd"
$ws.Range("C5").Value2 = "ASX"
$ws.Range("D5").Value2 = "2022-10-23T04:00:00.000Z"
$ws.Range("E5").Value2 = 0
$ws.Range("F5").Value2 = "2022-12-06T05:00:00.000Z"
$ws.Range("G5").Value2 = 0
$ws.Range("H5").Value2 = 0
$ws.Range("I5").Value2 = 4
$ws.Range("J5").Value2 = "AUD"
$ws.Range("K5").Value2 = 5000000
$ws.Range("L5").Value2 = "Amended"

# Re-enter the array formula so its `ref` shrinks from A5:L6 to A5:L5.
$ws.Range("A5:L5").FormulaArray = "=_xldudf_STREAM_PARSEJSONRANGE(A3)"

# Old row 6 (previously the ACE.AX record) is no longer part of the spill -
# clear it out entirely.
$ws.Range("A6:L6").ClearContents()

# --- Row 10 formula -------------------------------------------------------
# Add the second TRUE argument to the STREAM_PARSEJSONRANGE call that feeds
# rows 10:12; keep the anchor/ref the same (A10:L12).
$ws.Range("A10:L12").FormulaArray = '=_xldudf_STREAM_PARSEJSONRANGE(_xldudf_STREAM_EXTRACT(A1, "data"), TRUE)'

# --- Selection --------------------------------------------------------
$ws.Range("B16").Select()
